$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows scraped from MV -datos- (05-08-2021 .. 03-09-2021)
$data = @(
    @("05-08-2021", 3598, 6788),
    @("06-08-2021", 3486, 6803),
    @("09-08-2021", 4128, 7796),
    @("10-08-2021", 3487, 7895),
    @("11-08-2021", 3576, 7671),
    @("12-08-2021", 3571, 8804),
    @("13-08-2021", 3786, 7761),
    @("16-08-2021", 3469, 7056),
    @("17-08-2021", 3461, 8757),
    @("18-08-2021", 3404, 7149),
    @("19-08-2021", 3666, 8050),
    @("20-08-2021", 3494, 6994),
    @("23-08-2021", 3354, 5558),
    @("24-08-2021", 3375, 8276),
    @("25-08-2021", 3660, 7815),
    @("26-08-2021", 3623, 6811),
    @("27-08-2021", 3315, 6614),
    @("30-08-2021", 3515, 6001),
    @("31-08-2021", 3449, 7765),
    @("01-09-2021", 3809, 10926),
    @("02-09-2021", 3766, 9841),
    @("03-09-2021", 4063, 8005)
)

$startRow = 150
$endRow = $startRow + $data.Count - 1

# Format column A as Text first so the dd-mm-yyyy strings are not
# auto-converted into date serial numbers by Excel, then restore the
# default (Normal) style so the new cells match the rest of the sheet.
$dateRange = $ws.Range("A" + $startRow + ":A" + $endRow)
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$dateRange.Style = "Normal"

"Added rows " + $startRow + "-" + $endRow
